# Populate the Active / InActive / InComplete sheets with the real product
# rows (sku, item name, asin, product id, relation, avg landed cost).
# Every value must land as literal text (matching the source workbook's
# shared-string cells), so each cell is forced to the "@" (Text) number
# format before the value is assigned - otherwise Excel auto-converts
# numeric-looking strings (UPCs, landed-cost figures, ...) into numbers.
# The style is reset back to "Normal" right after so the cell is left with
# no lingering explicit style, matching plain data-entry cells elsewhere in
# the sheet.

function Set-TextCell {
    param($ws, [string]$addr, [string]$val)
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

function Set-RowValues {
    param($ws, [int]$row, [string[]]$values)
    $cols = @("A", "B", "C", "D", "E", "F")
    for ($i = 0; $i -lt $values.Length; $i++) {
        $addr = $cols[$i] + $row
        $val = $values[$i]
        Set-TextCell $ws $addr $val
    }
}

$wb = $excel.ActiveWorkbook

$active = $wb.Worksheets.Item("Active")
$inactive = $wb.Worksheets.Item("InActive")
$incomplete = $wb.Worksheets.Item("InComplete")

# --- Active sheet --------------------------------------------------------
Set-RowValues $active 2 @(
    "1B-YHD0-06JY",
    "Baby Diaper Caddy Organizer – Nursery Basket with Convenient Leather Handles, Storage Bin – Durable, Portable Changing Table Diaper Storage + Bonus Insulated Wipe Carrier by Cartik™",
    "B07C27CG18",
    "661708972459",
    "missing",
    "4.6"
)
Set-RowValues $active 3 @(
    "EL-NTZN-UZAR",
    "Baby Diaper Caddy Organizer – Nursery Basket with Convenient Leather Handles – Durable, Portable Changing Table Diaper Storage (2 Pack)",
    "B09PGLKMK9",
    "B09PGLKMK9",
    "child",
    "8"
)
Set-RowValues $active 4 @(
    "HZ-QVQR-JS29",
    "Cartik 2 Pack Backseat Car Organizer for Kids, Babies and Toddlers, with Tablet Holder by iPad Touch Screen, Fit to Baby Stroller, Large Storage, Kick Mat, Back Seat Protector, Organizer eBook",
    "B07GNRHN2Q",
    "661708972442",
    "child",
    "7.2"
)
Set-RowValues $active 5 @(
    "LD-MT1T-ZNZU",
    "Cartik™ Backseat Car Organizer Kids, Babies Toddlers Tablet Holder iPad Touch Screen, Fit to Baby Stroller, Large Storage, Kick Mat, Back Seat Protector, Organizer eBook (one Pack)",
    "B076ZJX4SX",
    "B076ZJX4SX",
    "child",
    "3.6"
)

# --- InActive sheet -------------------------------------------------------
Set-RowValues $inactive 2 @(
    "5L-Y8DM-ULLO",
    "Diaper Caddy Organizer (old Diaper Caddy Organizer)",
    "B07F1X69HS",
    "661708972466",
    "child",
    "3"
)
Set-RowValues $inactive 3 @(
    "JB-GIX4-MKM5",
    "Diaper Caddy Organizer (old Diaper Caddy Organizer)",
    "B07F1X69HS",
    "B07F1X69HS",
    "child",
    "3"
)
Set-RowValues $inactive 4 @(
    "O2-WSWS-RNP8",
    "Cartik Backseat Car Organizer for Kids, Babies and Toddlers, with Tablet Holder by iPad Touch Screen, Fit to Baby Stroller, Large Storage, Kick Mat, Back Seat Protector, Organizer eBook",
    "B07FZQRZZF",
    "B07FZQRZZF",
    "parent",
    "6"
)

# --- InComplete sheet -------------------------------------------------------
Set-RowValues $incomplete 2 @(
    "U1-4AA3-M779",
    "Diaper Caddy Organizer",
    "B09YFP28G8",
    "B09YFP28G8",
    "parent",
    "3"
)
